$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper-free, explicit per-cell updates to mirror the authoritative diff.
# D-column price cells are forced to Text before/after the write so that
# numeric-looking strings (e.g. "335.05") are stored as shared-string text
# (matching the source t="inlineStr" cells) instead of being coerced to
# numbers by Excel's automatic type detection. ClearFormats() afterwards
# resets the style index back to the sheet default (0) so no incidental
# style/number-format change is introduced.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# Row 2
Set-TextValue $ws.Range('D2') '30.435.17'
$ws.Range('E2').Value = '  +0.70%  '

# Row 3
Set-TextValue $ws.Range('D3') '2.111.26'
$ws.Range('E3').Value = '  +1.92%  '

# Row 4
$ws.Range('E4').Value = '  +0.44%  '

# Row 5
Set-TextValue $ws.Range('D5') '335.05'
$ws.Range('E5').Value = '  +2.78%  '

# Row 6
$ws.Range('E6').Value = '  +0.47%  '

# Row 7
Set-TextValue $ws.Range('D7') '0.5237'
$ws.Range('E7').Value = '  +1.28%  '

# Row 8
Set-TextValue $ws.Range('D8') '0.4544'
$ws.Range('E8').Value = '  +5.34%  '

# Row 9
Set-TextValue $ws.Range('D9') '53.51'
$ws.Range('E9').Value = '  +17.38%  '

# Row 10
Set-TextValue $ws.Range('D10') '0.08927'
$ws.Range('E10').Value = '  +2.33%  '

# Row 11
$ws.Range('E11').Value = '  +2.53%  '

# Row 12
Set-TextValue $ws.Range('D12') '24.40'
$ws.Range('E12').Value = '  +1.32%  '

# Row 13
Set-TextValue $ws.Range('D13') '2.101.25'
$ws.Range('E13').Value = '  +1.57%  '

# Row 14
Set-TextValue $ws.Range('D14') '6.849'
$ws.Range('E14').Value = '  +3.48%  '

# Row 15
Set-TextValue $ws.Range('D15') '8.092'
$ws.Range('E15').Value = '  +6.13%  '

# Row 16
Set-TextValue $ws.Range('D16') '96.78'
$ws.Range('E16').Value = '  +2.23%  '

# Row 17
$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue $ws.Range('D17') '0.00001148'
$ws.Range('E17').Value = '  +3.10%  '

# Row 18
$ws.Range('B18').Value = 'BinanceUSD'
$ws.Range('C18').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextValue $ws.Range('D18') '1.006'
$ws.Range('E18').Value = '  +0.52%  '

# Row 19
Set-TextValue $ws.Range('D19') '0.06651'
$ws.Range('E19').Value = '  +0.84%  '

# Row 20
$ws.Range('E20').Value = '  +3.36%  '

# Row 21
$ws.Range('E21').Value = '  +0.16%  '

# Row 22
Set-TextValue $ws.Range('D22') '6.367'
$ws.Range('E22').Value = '  +2.61%  '

# Row 23
Set-TextValue $ws.Range('D23') '30.491.47'
$ws.Range('E23').Value = '  +0.75%  '

# Row 24
$ws.Range('E24').Value = '  +1.78%  '

# Row 25
Set-TextValue $ws.Range('D25') '2.366'
$ws.Range('E25').Value = '  +4.49%  '

# Row 26
Set-TextValue $ws.Range('D26') '2.344.02'
$ws.Range('E26').Value = '  +1.78%  '

# Row 27
Set-TextValue $ws.Range('D27') '22.44'
$ws.Range('E27').Value = '  +1.72%  '

# Row 28
Set-TextValue $ws.Range('D28') '2.578'
$ws.Range('E28').Value = '  +3.61%  '

# Row 29
Set-TextValue $ws.Range('D29') '163.69'
$ws.Range('E29').Value = '  +1.29%  '

# Row 30
Set-TextValue $ws.Range('D30') '133.85'
$ws.Range('E30').Value = '  +2.84%  '

# Row 31
Set-TextValue $ws.Range('D31') '1.241'
$ws.Range('E31').Value = '  +5.19%  '

# Row 32
Set-TextValue $ws.Range('D32') '0.1075'
$ws.Range('E32').Value = '  +1.30%  '

# Row 33
Set-TextValue $ws.Range('D33') '1.696'
$ws.Range('E33').Value = '  +12.89%  '

# Row 34
Set-TextValue $ws.Range('D34') '6.344'
$ws.Range('E34').Value = '  +4.90%  '

# Row 35
Set-TextValue $ws.Range('D35') '3.936'
$ws.Range('E35').Value = '  +2.60%  '

# Row 36
Set-TextValue $ws.Range('D36') '10.51'
$ws.Range('E36').Value = '  +9.99%  '

# Row 37
Set-TextValue $ws.Range('D37') '0.02597'
$ws.Range('E37').Value = '  +1.86%  '

# Row 38
Set-TextValue $ws.Range('D38') '5.663'
$ws.Range('E38').Value = '  +5.04%  '

# Row 39
Set-TextValue $ws.Range('D39') '0.06838'
$ws.Range('E39').Value = '  +4.26%  '

# Row 40
Set-TextValue $ws.Range('D40') '0.2306'
$ws.Range('E40').Value = '  +4.09%  '

# Row 41
Set-TextValue $ws.Range('D41') '12.76'
$ws.Range('E41').Value = '  +2.48%  '

# Row 42
Set-TextValue $ws.Range('D42') '0.6919'
$ws.Range('E42').Value = '  +4.43%  '

# Row 43
Set-TextValue $ws.Range('D43') '1.251'
$ws.Range('E43').Value = '  +1.42%  '

# Row 44
Set-TextValue $ws.Range('D44') '2.356'
$ws.Range('E44').Value = '  +8.30%  '

# Row 45
Set-TextValue $ws.Range('D45') '1.004'
$ws.Range('E45').Value = '  +0.43%  '

# Row 46
Set-TextValue $ws.Range('D46') '14.13'
$ws.Range('E46').Value = '  +1.54%  '

# Row 47
Set-TextValue $ws.Range('D47') '0.6406'
$ws.Range('E47').Value = '  +2.26%  '

# Row 48
Set-TextValue $ws.Range('D48') '3.669'
$ws.Range('E48').Value = '  +2.08%  '

# Row 49
$ws.Range('B49').Value = 'EOS'
$ws.Range('C49').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
Set-TextValue $ws.Range('D49') '1.253'
$ws.Range('E49').Value = '  +1.95%  '

# Row 50
$ws.Range('B50').Value = 'WOONetwork'
$ws.Range('C50').Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
Set-TextValue $ws.Range('D50') '0.3457'
$ws.Range('E50').Value = '  +27.14%  '

# Row 51
$ws.Range('B51').Value = 'BabyDogeCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextValue $ws.Range('D51') '0.00000000344'
$ws.Range('E51').Value = '  +20.88%  '
